$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to swap between paired rows (everything except column A, which holds
# the row's sequential index, and C/D which are identical within each pair anyway).
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$pairs = @(
    @(60, 61),
    @(77, 78),
    @(132, 133),
    @(140, 141),
    @(186, 187),
    @(243, 244)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
